$d = $word.ActiveDocument

$replacements = @(
    @{old='45×77=3465'; new='58×40=2320'},
    @{old='97×99=9603'; new='71×20=1420'},
    @{old='60×38=2280'; new='78×42=3276'},
    @{old='15×62=930'; new='40×47=1880'},
    @{old='25×85=2125'; new='54×18=972'},
    @{old='81×58=4698'; new='33×67=2211'},
    @{old='37×18=666'; new='70×12=840'},
    @{old='53×88=4664'; new='28×17=476'},
    @{old='44×31=1364'; new='79×66=5214'},
    @{old='38×43=1634'; new='64×59=3776'},
    @{old='91×38=3458'; new='97×56=5432'},
    @{old='27×58=1566'; new='88×25=2200'},
    @{old='29×73=2117'; new='30×15=450'},
    @{old='90×50=4500'; new='70×30=2100'},
    @{old='79×85=6715'; new='20×30=600'},
    @{old='46×38=1748'; new='71×49=3479'},
    @{old='42×47=1974'; new='16×81=1296'},
    @{old='13×73=949'; new='54×38=2052'},
    @{old='76×90=6840'; new='84×82=6888'},
    @{old='93×31=2883'; new='96×23=2208'},
    @{old='70×29=2030'; new='16×50=800'},
    @{old='94×11=1034'; new='58×14=812'},
    @{old='45×89=4005'; new='62×47=2914'},
    @{old='65×29=1885'; new='26×90=2340'},
    @{old='40×56=2240'; new='34×26=884'}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
